# Swap the TORENBEEK_1982 / SFORZA comparison rows on the FUSELAGE and WING
# sheets of the Balance workbook ("Some work on time constant ATR72_ATsi").
#
# On FUSELAGE (rows 23-24) and WING (rows 23-24, and rows 27-28) the two
# rows' values are exchanged while the row labels (A column) stay put -
# i.e. the TORENBEEK_1982 figure moves into the row that used to hold the
# SFORZA figure, and vice versa.

$wb = $excel.ActiveWorkbook

function Swap-CellValues($ws, $cellA, $cellB) {
    $rangeA = $ws.Range($cellA)
    $rangeB = $ws.Range($cellB)
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# FUSELAGE sheet: rows 23 and 24 (column C)
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
Swap-CellValues $wsFuselage "C23" "C24"

# WING sheet: rows 23/24 and rows 27/28 (column C)
$wsWing = $wb.Worksheets.Item("WING")
Swap-CellValues $wsWing "C23" "C24"
Swap-CellValues $wsWing "C27" "C28"
